$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("supervision")
if (-not $ws) { $ws = $wb.ActiveSheet }

# Add new row (row 12) for Adele Chu
$ws.Range("A12").Value = "Adele Chu"

$ws.Range("F12").NumberFormat = "mm-dd-yy"
$ws.Range("F12").Value = 45658

$ws.Range("G12").NumberFormat = "m/d/yy h:mm"
$ws.Range("G12").Formula = "=NOW()"

$ws.Range("J12").Value = "Cosupervision"

# Update selection / zoom to match the edited view
$ws.Range("K12").Select()
$excel.ActiveWindow.Zoom = 187

$wb.Save()
